# Update market price data (currentAveragePrice / profit columns) across leve-crafting sheets
# per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 8338691
$ws.Range("I32").Value = 7665
$ws.Range("K32").Value = 7665
$ws.Range("M32").Value = -7339
$ws.Range("H70").Value = 3232.8076
$ws.Range("I70").Value = 2914.0833
$ws.Range("K70").Value = 8742.249899999999
$ws.Range("M70").Value = -8472.249899999999
$ws.Range("H73").Value = 3232.8076
$ws.Range("I73").Value = 2914.0833
$ws.Range("K73").Value = 8742.249899999999
$ws.Range("M73").Value = -7806.249899999999
$ws.Range("H76").Value = 4612.846
$ws.Range("I76").Value = 4691.579
$ws.Range("K76").Value = 4691.579
$ws.Range("M76").Value = -4376.579
$ws.Range("H79").Value = 4612.846
$ws.Range("I79").Value = 4691.579
$ws.Range("K79").Value = 4691.579
$ws.Range("M79").Value = -3599.579
$ws.Range("H99").Value = 71207.92999999999
$ws.Range("I99").Value = 559.9091
$ws.Range("J99").Value = 265490
$ws.Range("K99").Value = 1679.7273
$ws.Range("L99").Value = 796470
$ws.Range("M99").Value = -181.7273
$ws.Range("N99").Value = -799466
$ws.Range("H100").Value = 4071.5715
$ws.Range("I100").Value = 2052.8
$ws.Range("K100").Value = 2052.8
$ws.Range("M100").Value = -1511.8
$ws.Range("H101").Value = 2467
$ws.Range("I101").Value = 960.4
$ws.Range("K101").Value = 2881.2
$ws.Range("M101").Value = -1259.2
$ws.Range("H107").Value = 371.25
$ws.Range("I107").Value = 315
$ws.Range("K107").Value = 315
$ws.Range("M107").Value = 1605
$ws.Range("H112").Value = 1975.4546
$ws.Range("J112").Value = 1988.8372
$ws.Range("L112").Value = 5966.5116
$ws.Range("N112").Value = -8182.5116
$ws.Range("H135").Value = 3589.65
$ws.Range("I135").Value = 3934.111
$ws.Range("K135").Value = 35406.999
$ws.Range("M135").Value = -32871.999
$ws.Range("H137").Value = 3191.9656
$ws.Range("I137").Value = 2760.7083
$ws.Range("K137").Value = 8282.124899999999
$ws.Range("M137").Value = -5732.124899999999
$ws.Range("H138").Value = 2983.7793
$ws.Range("I138").Value = 1435.1305
$ws.Range("J138").Value = 3775.311
$ws.Range("K138").Value = 4305.3915
$ws.Range("L138").Value = 11325.933
$ws.Range("M138").Value = 834.6085000000003
$ws.Range("N138").Value = -21605.933

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9536.589
$ws.Range("I32").Value = 8930.014999999999
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 8930.014999999999
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -8643.014999999999
$ws.Range("N32").Value = -20574
$ws.Range("H63").Value = 3663
$ws.Range("I63").Value = 3829.1428
$ws.Range("K63").Value = 3829.1428
$ws.Range("M63").Value = -3143.1428
$ws.Range("H66").Value = 3663
$ws.Range("I66").Value = 3829.1428
$ws.Range("K66").Value = 19145.714
$ws.Range("M66").Value = -15713.714
$ws.Range("H97").Value = 1227.4286
$ws.Range("I97").Value = 1308.6842
$ws.Range("K97").Value = 1308.6842
$ws.Range("M97").Value = -812.6841999999999
$ws.Range("H132").Value = 3345.4902
$ws.Range("I132").Value = 3023.721
$ws.Range("J132").Value = 5075
$ws.Range("K132").Value = 9071.163
$ws.Range("L132").Value = 15225
$ws.Range("M132").Value = -6541.163
$ws.Range("N132").Value = -20285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2622.125
$ws.Range("I22").Value = 4121
$ws.Range("J22").Value = 124
$ws.Range("K22").Value = 4121
$ws.Range("L22").Value = 124
$ws.Range("M22").Value = -3948
$ws.Range("N22").Value = -470
$ws.Range("H44").Value = 950
$ws.Range("I44").Value = 950
$ws.Range("K44").Value = 950
$ws.Range("M44").Value = -453
$ws.Range("H86").Value = 3877.476
$ws.Range("I86").Value = 2188
$ws.Range("J86").Value = 6130.1113
$ws.Range("K86").Value = 2188
$ws.Range("L86").Value = 6130.1113
$ws.Range("M86").Value = -1065
$ws.Range("N86").Value = -8376.1113
$ws.Range("H89").Value = 3877.476
$ws.Range("I89").Value = 2188
$ws.Range("J89").Value = 6130.1113
$ws.Range("K89").Value = 10940
$ws.Range("L89").Value = 30650.5565
$ws.Range("M89").Value = -5324
$ws.Range("N89").Value = -41882.5565
$ws.Range("H94").Value = 7119.9
$ws.Range("I94").Value = 6525
$ws.Range("K94").Value = 6525
$ws.Range("M94").Value = -6074
$ws.Range("H96").Value = 98999
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H99").Value = 25371.059
$ws.Range("I99").Value = 26644.25
$ws.Range("K99").Value = 26644.25
$ws.Range("M99").Value = -25146.25
$ws.Range("H105").Value = 3204.4443
$ws.Range("I105").Value = 2283.3076
$ws.Range("J105").Value = 5599.4
$ws.Range("K105").Value = 2283.3076
$ws.Range("L105").Value = 5599.4
$ws.Range("M105").Value = -536.3076000000001
$ws.Range("N105").Value = -9093.4
$ws.Range("H134").Value = 3664.0322
$ws.Range("I134").Value = 2453.5417
$ws.Range("K134").Value = 7360.625100000001
$ws.Range("M134").Value = -4825.625100000001
$ws.Range("H137").Value = 69727.734
$ws.Range("J137").Value = 69727.734
$ws.Range("L137").Value = 69727.734
$ws.Range("N137").Value = -79927.734

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3158.853
$ws.Range("I132").Value = 2599.2415
$ws.Range("J132").Value = 6404.6
$ws.Range("K132").Value = 7797.7245
$ws.Range("L132").Value = 19213.8
$ws.Range("M132").Value = -5267.7245
$ws.Range("N132").Value = -24273.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 157.8125
$ws.Range("I12").Value = 116.333336
$ws.Range("J12").Value = 174.04347
$ws.Range("K12").Value = 349.000008
$ws.Range("L12").Value = 522.13041
$ws.Range("M12").Value = -176.000008
$ws.Range("N12").Value = -868.13041
$ws.Range("H39").Value = 3526.7646
$ws.Range("J39").Value = 3937
$ws.Range("L39").Value = 11811
$ws.Range("N39").Value = -12399
$ws.Range("H119").Value = 6666.6665
$ws.Range("I119").Value = 2000
$ws.Range("K119").Value = 6000
$ws.Range("M119").Value = -1162
$ws.Range("H141").Value = 2900
$ws.Range("I141").Value = 2900
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8700
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = -3520
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 15496.5
$ws.Range("J38").Value = 15496.5
$ws.Range("L38").Value = 15496.5
$ws.Range("N38").Value = -16422.5
$ws.Range("H134").Value = 73705.75
$ws.Range("J134").Value = 73705.75
$ws.Range("L134").Value = 221117.25
$ws.Range("N134").Value = -226187.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2790.1292
$ws.Range("I7").Value = 1499.4762
$ws.Range("J7").Value = 5500.5
$ws.Range("K7").Value = 1499.4762
$ws.Range("L7").Value = 5500.5
$ws.Range("M7").Value = -1387.4762
$ws.Range("N7").Value = -5724.5
$ws.Range("H22").Value = 1270.8334
$ws.Range("I22").Value = 1025
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 1025
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -730
$ws.Range("N22").Value = -3090
$ws.Range("H27").Value = 1270.8334
$ws.Range("I27").Value = 1025
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 1025
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -918
$ws.Range("N27").Value = -2714
$ws.Range("H82").Value = 2563.7188
$ws.Range("J82").Value = 4394.1
$ws.Range("L82").Value = 4394.1
$ws.Range("N82").Value = -5116.1
$ws.Range("H85").Value = 2563.7188
$ws.Range("J85").Value = 4394.1
$ws.Range("L85").Value = 4394.1
$ws.Range("N85").Value = -6890.1
$ws.Range("H93").Value = 295916.66
$ws.Range("I93").Value = 1698.8462
$ws.Range("K93").Value = 1698.8462
$ws.Range("M93").Value = -450.8462
$ws.Range("H100").Value = 53753.184
$ws.Range("I100").Value = 102489.63
$ws.Range("J100").Value = 5016.727
$ws.Range("K100").Value = 102489.63
$ws.Range("L100").Value = 5016.727
$ws.Range("M100").Value = -101948.63
$ws.Range("N100").Value = -6098.727
$ws.Range("H126").Value = 2790.1292
$ws.Range("I126").Value = 1499.4762
$ws.Range("J126").Value = 5500.5
$ws.Range("K126").Value = 4498.4286
$ws.Range("L126").Value = 16501.5
$ws.Range("M126").Value = -2028.4286
$ws.Range("N126").Value = -21441.5
$ws.Range("H136").Value = 4545.763
$ws.Range("I136").Value = 2320.0527
$ws.Range("K136").Value = 6960.158100000001
$ws.Range("M136").Value = -4410.158100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3871.5
$ws.Range("I62").Value = 3666.3333
$ws.Range("J62").Value = 4076.6667
$ws.Range("K62").Value = 3666.3333
$ws.Range("L62").Value = 4076.6667
$ws.Range("M62").Value = -3042.3333
$ws.Range("N62").Value = -5324.6667
$ws.Range("H65").Value = 3871.5
$ws.Range("I65").Value = 3666.3333
$ws.Range("J65").Value = 4076.6667
$ws.Range("K65").Value = 18331.6665
$ws.Range("L65").Value = 20383.3335
$ws.Range("M65").Value = -15211.6665
$ws.Range("N65").Value = -26623.3335
$ws.Range("H132").Value = 2325.5518
$ws.Range("I132").Value = 2159.22
$ws.Range("J132").Value = 3365.125
$ws.Range("K132").Value = 6477.66
$ws.Range("L132").Value = 10095.375
$ws.Range("M132").Value = -3947.66
$ws.Range("N132").Value = -15155.375

